$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) / Volume(1h) (E) columns, and swap the
# EthereumClassic / Stellar rows (29 <-> 30) per the latest scrape.
#
# Numeric-looking price strings (e.g. '0.694') must stay TEXT, matching the source
# feed's formatting (it also stores thousand-dotted prices like '35.278.42' as
# text) -- so for those cells we flip the cell to text format before writing the
# value, then restore the Normal style so no stray formatting is left behind.

$ws.Range('D2').Value = '35.278.42'
$ws.Range('D3').Value = '1.899.10'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.694'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +9.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '245.14'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '40.69'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.348'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.43%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '52.98'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +10.94%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0721'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '12.42'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('D16').Value = '1.907.62'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '35.304.26'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.15'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '0.0₃0818'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '240.51'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('E23').Value = '  -2.74%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.32'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.62%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '168.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.61'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.130'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.69%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '18.31'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('D31').Value = '4.140.38'
$ws.Range('E31').Value = '  +21.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.14'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.33%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0568'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.924'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.82'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.50%  '
$ws.Range('E38').Value = '  +11.04%  '
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0665'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +10.79%  '
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '15.98'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '89.54'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('D45').Value = '1.350.13'
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.73'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.42'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '45.96'
$ws.Range('D50').Style = "Normal"
$ws.Range('E51').Value = '  -2.55%  '
